$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the rows that were removed in the "shortened" action list.
# Delete bottom-up so row numbers of earlier rows stay valid.

# Removed: Pinch, thumb down (High), Pinch, thumb down (No contact),
#          Pinch, thumb up (Low/High/No contact), Draw word, Pinch-zoom
$ws.Range("A30:A36").EntireRow.Delete()

# Removed: Press fingers (Low), Press fingers (High)
$ws.Range("A16:A17").EntireRow.Delete()

# Removed: Index finger, pull towards / push left / push away
$ws.Range("A6:A8").EntireRow.Delete()

# A few of the remaining rows need their (explicit / custom) row height
# increased now that the sheet has fewer rows.
$ws.Rows.Item(6).RowHeight = 27.75
$ws.Rows.Item(10).RowHeight = 27.75
$ws.Rows.Item(11).RowHeight = 27.75
$ws.Rows.Item(12).RowHeight = 27.75
